# Daily attendance processing - normalize "Recorded By" (column G) ordering.
# For every data row, the comma-separated list of recorders in column G is
# reversed, unless the list already starts with the exact token "System"
# (case-sensitive) - in which case it is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $used.Row + $used.Rows.Count - 1

# Column G is the "Recorded By" column (7th column).
$col = 7

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $raw = $cell.Value2

    if ($raw -eq $null) {
        continue
    }

    $text = [string]$raw

    if ($text -notlike "*,*") {
        continue
    }

    $parts = $text.Split(",")
    $n = $parts.Count

    for ($i = 0; $i -lt $n; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts[0].Equals("System")) {
        continue
    }

    $reversed = @()
    for ($i = $n - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }

    $newText = [string]::Join(", ", $reversed)
    $cell.Value2 = $newText
}
